# Auto-generated edit script: update cryptos list values per commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain-decimal price cells to remain Text (matches original inlineStr type)
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Apply updated cell values
$ws.Range("D2").Value = '68.205.98'
$ws.Range("E2").Value = '  +1.38%  '
$ws.Range("D3").Value = '3.904.55'
$ws.Range("E3").Value = '  +0.76%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '479.08'
$ws.Range("E5").Value = '  +1.65%  '
$ws.Range("D6").Value = '144.07'
$ws.Range("E6").Value = '  -0.96%  '
$ws.Range("D7").Value = '0.619'
$ws.Range("E7").Value = '  -2.62%  '
$ws.Range("E8").Value = '  -0.19%  '
$ws.Range("D9").Value = '0.720'
$ws.Range("E9").Value = '  -3.81%  '
$ws.Range("D10").Value = '0.166'
$ws.Range("E10").Value = '  +7.34%  '
$ws.Range("D11").Value = '0.0000349'
$ws.Range("E11").Value = '  +11.63%  '
$ws.Range("D12").Value = '42.36'
$ws.Range("E12").Value = '  -2.99%  '
$ws.Range("D13").Value = '10.44'
$ws.Range("E13").Value = '  -0.07%  '
$ws.Range("D14").Value = '4.530.75'
$ws.Range("E14").Value = '  +0.73%  '
$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '3.930.48'
$ws.Range("E15").Value = '  +1.61%  '
$ws.Range("B16").Value = 'Uniswap'
$ws.Range("C16").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D16").Value = '14.55'
$ws.Range("E16").Value = '  -1.87%  '
$ws.Range("E17").Value = '  -0.37%  '
$ws.Range("E18").Value = '  -2.40%  '
$ws.Range("D19").Value = '1.12'
$ws.Range("E19").Value = '  -3.76%  '
$ws.Range("D20").Value = '68.252.59'
$ws.Range("E20").Value = '  +1.19%  '
$ws.Range("D21").Value = '431.48'
$ws.Range("E21").Value = '  -1.12%  '
$ws.Range("D22").Value = '14.56'
$ws.Range("E22").Value = '  -2.53%  '
$ws.Range("D23").Value = '3.35'
$ws.Range("E23").Value = '  +1.26%  '
$ws.Range("D24").Value = '87.05'
$ws.Range("E24").Value = '  -2.59%  '
$ws.Range("D25").Value = '11.64'
$ws.Range("E25").Value = '  +16.46%  '
$ws.Range("D26").Value = '3.55'
$ws.Range("E26").Value = '  -1.55%  '
$ws.Range("D27").Value = '37.94'
$ws.Range("E27").Value = '  -0.25%  '
$ws.Range("D28").Value = '10.17'
$ws.Range("E28").Value = '  -0.26%  '
$ws.Range("D29").Value = '5.82'
$ws.Range("E29").Value = '  +6.05%  '
$ws.Range("D30").Value = '699.85'
$ws.Range("E30").Value = '  -4.67%  '
$ws.Range("E31").Value = '  -3.79%  '
$ws.Range("D32").Value = '13.21'
$ws.Range("E32").Value = '  -4.90%  '
$ws.Range("E33").Value = '  +3.34%  '
$ws.Range("D34").Value = '0.0₃0903'
$ws.Range("E34").Value = '  +30.78%  '
$ws.Range("D35").Value = '41.14'
$ws.Range("E35").Value = '  -8.10%  '
$ws.Range("D36").Value = '59.15'
$ws.Range("E36").Value = '  +1.85%  '
$ws.Range("B37").Value = 'NEARProtocol'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D37").Value = '5.65'
$ws.Range("E37").Value = '  +2.72%  '
$ws.Range("B38").Value = 'Kaspa'
$ws.Range("C38").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D38").Value = '0.150'
$ws.Range("E38").Value = '  -8.48%  '
$ws.Range("D39").Value = '0.997'
$ws.Range("E39").Value = '  -0.30%  '
$ws.Range("D40").Value = '0.0471'
$ws.Range("E40").Value = '  -2.99%  '
$ws.Range("D41").Value = '3.06'
$ws.Range("E41").Value = '  +11.26%  '
$ws.Range("D42").Value = '2.74'
$ws.Range("E42").Value = '  +6.66%  '
$ws.Range("E43").Value = '  +2.18%  '
$ws.Range("E44").Value = '  -0.49%  '
$ws.Range("E45").Value = '  -3.34%  '
$ws.Range("E46").Value = '  +0.03%  '
$ws.Range("D47").Value = '3.41'
$ws.Range("E47").Value = '  -1.80%  '
$ws.Range("D48").Value = '2.13'
$ws.Range("E48").Value = '  -1.40%  '
$ws.Range("D49").Value = '146.56'
$ws.Range("E49").Value = '  +1.57%  '
$ws.Range("D50").Value = '3.13'
$ws.Range("E50").Value = '  -5.03%  '
$ws.Range("E51").Value = '  -3.16%  '
